$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = "L_ERSTT_12"
$ws.Range("D4").Value = "L_BMEL_3"

# Row 19
$ws.Range("C19").Value = "L_BMZ_2"
$ws.Range("D19").Value = "Q_GIZ"
$ws.Range("E19").Value = "Q_CEVAL"

# Row 26
$ws.Range("C26").Value = "L_DSTTS_27"

# Row 31
$ws.Range("C31").Value = "Q_DESTATIS"

# Row 32
$ws.Range("C32").Value = "L_GIZ_1"
$ws.Range("D32").Value = "Q_GIZ"

# Row 38
$ws.Range("C38").Value = "L_RDB_1"
$ws.Range("D38").Value = "L_RDB_2"
$ws.Range("E38").Value = "L_RDB_3"
$ws.Range("F38").Value = "Q_BBSR"
$ws.Range("G38").Value = "Q_JHT"

# Row 39
$ws.Range("C39").Value = "L_RDB_1"
$ws.Range("D39").Value = "L_RDB_2"
$ws.Range("E39").Value = "L_RDB_3"
$ws.Range("F39").Value = "Q_BBSR"
$ws.Range("G39").Value = "Q_JHT"

# Row 40
$ws.Range("D40").Value = "L_IFEU_1"

# Row 41
$ws.Range("D41").Value = "L_IFEU_1"

# Row 43
$ws.Range("C43").Value = "L_DSTTS_24"
$ws.Range("D43").Value = "L_ERSTT_11"

# Row 45
$ws.Range("C45").Value = "L_UBA_4"
$ws.Range("D45").Value = ""
$ws.Range("E45").Value = ""
$ws.Range("F45").Value = ""
$ws.Range("G45").Value = ""
$ws.Range("H45").Value = ""

# Row 46
$ws.Range("C46").Value = "L_DSTTS_25"

# Row 53
$ws.Range("C53").Value = "L_BFN_1"
